# conductor_grid.xlsx - update GRID sheet to latest template version.
# Adds DXINCRE_LEFT/DXINCRE_RIGHT split (was a single DXINCRE) and a new
# MAXNOD parameter row; refreshes the ITYMSH note text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5 (ITYMSH): refreshed note text, taller row to fit it ---------
$ws.Range("D5").Value2 = "flag for mesh type: 0 = fixed and uniform; 1 = fixed refined; 3 = adapted with initial refinement; - 1 from file; in this case the z coordinates of the conductor components must be exactly the same for each conudctor component objets.                                                                     "
$ws.Rows.Item(5).RowHeight = 60

# --- Prepare two fresh rows (12 and 13) with the same look as row 11 ---
# (row 11 already carries the exact label/unit/type/note/value style
# pattern we need to reuse for the two new rows; only copy A:E so we
# do not drag formatting across the whole 16384-column row)
$ws.Range("A11:E11").Copy() | Out-Null
$ws.Range("A12:E12").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("A11:E11").Copy() | Out-Null
$ws.Range("A13:E13").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
# Row 13's value cell (MAXNOD) keeps the plain "-"/"integer" look (General
# number format, no scientific notation) rather than the float style, so
# borrow it from column B/C of an existing row instead of column E.
$ws.Range("B11").Copy() | Out-Null
$ws.Range("E13").PasteSpecial(-4122) | Out-Null       # xlPasteFormats
$excel.CutCopyMode = 0

# --- Row 11: DXINCRE -> DXINCRE_LEFT ------------------------------------
$ws.Range("A11").Value2 = "DXINCRE_LEFT"
$ws.Range("D11").Value2 = "size increase ratio for the spatial mesh, used for the region to the left of the refined region."
$ws.Range("E11").Value2 = 1.2
$ws.Rows.Item(11).RowHeight = 30

# --- Row 12 (new): DXINCRE_RIGHT ----------------------------------------
$ws.Range("A12").Value2 = "DXINCRE_RIGHT"
$ws.Range("B12").Value2 = "-"
$ws.Range("C12").Value2 = "float"
$ws.Range("D12").Value2 = "size increase ratio for the spatial mesh, used for the region to the right of the refined region."
$ws.Range("E12").Value2 = 1.2
$ws.Rows.Item(12).RowHeight = 30

# --- Row 13 (new): MAXNOD ------------------------------------------------
$ws.Range("A13").Value2 = "MAXNOD"
$ws.Range("B13").Value2 = "-"
$ws.Range("C13").Value2 = "integer"
$ws.Range("D13").Value2 = " maximum number of nodes for conductor spatial discretization"
$ws.Range("E13").Value2 = 10001

# --- Selection cursor, matching the saved view of the authored file ----
$ws.Range("D23").Select() | Out-Null
